$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("D2").ClearContents()
$ws.Range("B2").Value = 43.137025930401123
$ws.Range("C2").Value = 21.834834647764225

$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 42.084078618391473

$ws.Range("B1:E3").Select() | Out-Null
